$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.960.97"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").Value = "2.564.57"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.573"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.83%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("E13").Value = "  +7.16%  "

$ws.Range("D14").Value = "2.954.49"
$ws.Range("E14").Value = "  -2.22%  "

$ws.Range("D15").Value = "2.538.19"
$ws.Range("E15").Value = "  -3.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.881"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("D18").Value = "42.960.97"
$ws.Range("E18").Value = "  -1.53%  "

$ws.Range("D19").Value = "0.0₃0995"
$ws.Range("E19").Value = "  +2.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.65%  "

$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0799"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.19%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +31.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0310"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.50%  "

$ws.Range("D45").Value = "2.102.78"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.66%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.808.84"
$ws.Range("E50").Value = "  -2.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.86%  "
